$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.874.48"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "1.584.18"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.08"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("E7").Value = "  -1.71%  "

$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.10"
$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").Value = "1.803.79"
$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("D13").Value = "1.592.13"
$ws.Range("E13").Value = "  -1.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("E15").Value = "  -2.09%  "

$ws.Range("D16").Value = "25.876.82"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.08"
$ws.Range("E18").Value = "  -1.89%  "

$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.87"
$ws.Range("E20").Value = "  +1.62%  "

$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("E24").Value = "  +1.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.71"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("E27").Value = "  -1.18%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("E29").Value = "  -1.84%  "

$ws.Range("E30").Value = "  -4.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("E31").Value = "  +0.44%  "

$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("E34").Value = "  +1.48%  "

$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("D36").Value = "1.096.19"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.777"
$ws.Range("E41").Value = "  -4.26%  "

$ws.Range("E42").Value = "  +7.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "93.05"
$ws.Range("E43").Value = "  -3.70%  "

$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("D45").Value = "1.717.30"
$ws.Range("E45").Value = "  -1.69%  "

$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  -2.89%  "

$ws.Range("E47").Value = "  +2.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.09"
$ws.Range("E48").Value = "  -0.69%  "

$ws.Range("E49").Value = "  -1.20%  "

$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("E51").Value = "  -0.23%  "
